# Add "NA" values under the duplicate_image_filename column (column E)
# for data rows 2-21 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2:E21").Value = "NA"
